# Fixed #517 MCell starting with a MPragraph leaves an empty paragraph at
# the begining of the cell.
#
# Each cell of the table currently begins with an empty leading paragraph
# (no run, no text) immediately followed by the paragraph that actually
# carries the cell's content. Remove that superfluous leading empty
# paragraph from every cell of every table so the content paragraph
# becomes the cell's first (and only) paragraph.

$d = $word.ActiveDocument

foreach ($t in $d.Tables) {
    foreach ($row in $t.Rows) {
        foreach ($cell in $row.Cells) {
            $cellParagraphs = $cell.Range.Paragraphs
            while ($cellParagraphs.Count -gt 1 -and `
                   $cellParagraphs.Item(1).Range.Text.Trim() -eq "") {
                $cellParagraphs.Item(1).Range.Delete()
                $cellParagraphs = $cell.Range.Paragraphs
            }
        }
    }
}
